$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update existing AgTests (F) and AgPosit (G) values for rows with revised figures
$ws.Range("F670").Value = 52968
$ws.Range("G670").Value = 931
$ws.Range("F671").Value = 32815
$ws.Range("G671").Value = 634
$ws.Range("F672").Value = 29981
$ws.Range("G672").Value = 600
$ws.Range("F674").Value = 29003
$ws.Range("G674").Value = 720
$ws.Range("F677").Value = 56860
$ws.Range("G677").Value = 871
$ws.Range("F678").Value = 34036
$ws.Range("G678").Value = 554
$ws.Range("F679").Value = 29634
$ws.Range("G679").Value = 548
$ws.Range("F680").Value = 28629
$ws.Range("G680").Value = 575
$ws.Range("F681").Value = 26607
$ws.Range("G681").Value = 611
$ws.Range("F684").Value = 57899
$ws.Range("G684").Value = 1306
$ws.Range("F685").Value = 34667
$ws.Range("G685").Value = 1063
$ws.Range("F686").Value = 34626
$ws.Range("G686").Value = 1171
$ws.Range("F687").Value = 31673
$ws.Range("G687").Value = 1164
$ws.Range("F688").Value = 32354
$ws.Range("G688").Value = 1386
$ws.Range("F691").Value = 63084
$ws.Range("G691").Value = 2858
$ws.Range("F692").Value = 41871
$ws.Range("G692").Value = 2728
$ws.Range("F693").Value = 39807
$ws.Range("G693").Value = 2779
$ws.Range("F694").Value = 37833
$ws.Range("G694").Value = 2824
$ws.Range("F695").Value = 37502
$ws.Range("G695").Value = 3182
$ws.Range("F698").Value = 71395
$ws.Range("G698").Value = 5953
$ws.Range("F699").Value = 43674
$ws.Range("G699").Value = 4349
$ws.Range("F700").Value = 44012
$ws.Range("G700").Value = 4366
$ws.Range("F701").Value = 42082
$ws.Range("G701").Value = 3891
$ws.Range("F702").Value = 36562
$ws.Range("G702").Value = 3969
$ws.Range("F705").Value = 56711
$ws.Range("G705").Value = 6430
$ws.Range("F706").Value = 40912
$ws.Range("G706").Value = 4992
$ws.Range("F707").Value = 39096
$ws.Range("G707").Value = 4659
$ws.Range("F708").Value = 35833
$ws.Range("G708").Value = 4201
$ws.Range("F709").Value = 32617
$ws.Range("G709").Value = 4019
$ws.Range("F712").Value = 52005
$ws.Range("G712").Value = 6359
$ws.Range("F713").Value = 37550
$ws.Range("G713").Value = 4796
$ws.Range("F714").Value = 32754
$ws.Range("G714").Value = 4039
$ws.Range("F715").Value = 32067
$ws.Range("G715").Value = 3608
$ws.Range("F716").Value = 30003
$ws.Range("G716").Value = 3720
$ws.Range("F719").Value = 45310
$ws.Range("G719").Value = 5323
$ws.Range("F720").Value = 31467
$ws.Range("G720").Value = 3550
$ws.Range("F721").Value = 28146
$ws.Range("G721").Value = 3180
$ws.Range("F722").Value = 28164
$ws.Range("G722").Value = 2916
$ws.Range("F723").Value = 22999
$ws.Range("G723").Value = 2820
$ws.Range("F726").Value = 36470
$ws.Range("G726").Value = 4209
$ws.Range("F727").Value = 25409
$ws.Range("G727").Value = 2834
$ws.Range("F728").Value = 24953
$ws.Range("G728").Value = 2635
$ws.Range("F729").Value = 23492
$ws.Range("G729").Value = 2539
$ws.Range("F730").Value = 19773
$ws.Range("G730").Value = 2350
$ws.Range("F733").Value = 32320
$ws.Range("G733").Value = 3756
$ws.Range("F734").Value = 23342
$ws.Range("G734").Value = 2560
$ws.Range("F735").Value = 19521
$ws.Range("G735").Value = 2281
$ws.Range("F736").Value = 19816
$ws.Range("G736").Value = 2208
$ws.Range("F737").Value = 18721
$ws.Range("G737").Value = 2315
$ws.Range("F740").Value = 25292
$ws.Range("G740").Value = 2782
$ws.Range("F741").Value = 19155
$ws.Range("G741").Value = 1944
$ws.Range("F742").Value = 17508
$ws.Range("G742").Value = 1707
$ws.Range("F743").Value = 18247
$ws.Range("G743").Value = 1637
$ws.Range("F744").Value = 14858
$ws.Range("G744").Value = 1609
$ws.Range("F747").Value = 22797
$ws.Range("G747").Value = 2405
$ws.Range("F748").Value = 17070
$ws.Range("G748").Value = 1547
$ws.Range("F749").Value = 14960
$ws.Range("G749").Value = 1486
$ws.Range("F750").Value = 15200
$ws.Range("G750").Value = 1360
$ws.Range("F751").Value = 12664
$ws.Range("G751").Value = 1387
$ws.Range("F754").Value = 21428
$ws.Range("G754").Value = 1962
$ws.Range("F755").Value = 13850
$ws.Range("G755").Value = 1301
$ws.Range("F756").Value = 13836
$ws.Range("G756").Value = 1076
$ws.Range("F757").Value = 13649
$ws.Range("G757").Value = 1010
$ws.Range("F770").Value = 9121
$ws.Range("G770").Value = 428
$ws.Range("F771").Value = 9223
$ws.Range("G771").Value = 407
$ws.Range("F772").Value = 2623
$ws.Range("F776").Value = 14807
$ws.Range("G776").Value = 670
$ws.Range("F777").Value = 10409
$ws.Range("G777").Value = 450
$ws.Range("F778").Value = 8974
$ws.Range("G778").Value = 358
$ws.Range("F779").Value = 7269
$ws.Range("G779").Value = 300
$ws.Range("F780").Value = 2673
$ws.Range("G780").Value = 130
$ws.Range("F781").Value = 2707
$ws.Range("G781").Value = 144
$ws.Range("F782").Value = 10439
$ws.Range("G782").Value = 414

# Fill in missing F783/G783 (row existed but lacked these values)
$ws.Range("F783").Value = 7488
$ws.Range("G783").Value = 239

# Append new row 784 with full data, matching style of column A (date format)
$ws.Range("A784").Value = 44678
$ws.Range("B784").Value = 1779096
$ws.Range("C784").Value = 5489
$ws.Range("D784").Value = 1125
$ws.Range("E784").Value = 19879
$ws.Range("F784").Value = 5525
$ws.Range("G784").Value = 176

# Apply date style (same as column A cells above) to the new row A784 cell
$ws.Range("A784").NumberFormat = $ws.Range("A783").NumberFormat
